$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H28").Value = 372.14285
$ws.Range("I28").Value = 392.5
$ws.Range("K28").Value = 392.5
$ws.Range("M28").Value = 92.5
$ws.Range("H33").Value = 191.75
$ws.Range("I33").Value = 193
$ws.Range("J33").Value = 188
$ws.Range("K33").Value = 193
$ws.Range("L33").Value = 188
$ws.Range("M33").Value = 36
$ws.Range("N33").Value = -646
$ws.Range("H125").Value = 1900
$ws.Range("I125").Value = 1900
$ws.Range("K125").Value = 17100
$ws.Range("M125").Value = -14640
$ws.Range("H129").Value = 1586.3334
$ws.Range("J129").Value = 1694.3182
$ws.Range("L129").Value = 5082.9546
$ws.Range("N129").Value = -15082.9546
$ws.Range("H137").Value = 85157.836
$ws.Range("I137").Value = 1397.6666
$ws.Range("J137").Value = 168918
$ws.Range("K137").Value = 4192.9998
$ws.Range("L137").Value = 506754
$ws.Range("M137").Value = -1642.9998
$ws.Range("N137").Value = -511854
$ws.Range("H138").Value = 3040.1833
$ws.Range("J138").Value = 2984.7188
$ws.Range("L138").Value = 8954.1564
$ws.Range("N138").Value = -19234.1564
$ws.Range("H141").Value = 1868815.6
$ws.Range("I141").Value = 2546839.5
$ws.Range("J141").Value = 4249.75
$ws.Range("K141").Value = 7640518.5
$ws.Range("L141").Value = 12749.25
$ws.Range("M141").Value = -7635338.5
$ws.Range("N141").Value = -23109.25

$ws = $wb.Worksheets("ARM")
$ws.Range("H74").Value = 803.381
$ws.Range("I74").Value = 585.3684
$ws.Range("J74").Value = 2874.5
$ws.Range("K74").Value = 585.3684
$ws.Range("L74").Value = 2874.5
$ws.Range("M74").Value = 288.6316
$ws.Range("N74").Value = -4622.5
$ws.Range("H77").Value = 803.381
$ws.Range("I77").Value = 585.3684
$ws.Range("J77").Value = 2874.5
$ws.Range("K77").Value = 2926.842
$ws.Range("L77").Value = 14372.5
$ws.Range("M77").Value = 1441.158
$ws.Range("N77").Value = -23108.5
$ws.Range("H128").Value = 30000
$ws.Range("J128").Value = 30000
$ws.Range("L128").Value = 30000
$ws.Range("N128").Value = -39960
$ws.Range("H132").Value = 1866.662
$ws.Range("I132").Value = 1486.2325
$ws.Range("J132").Value = 2450.8928
$ws.Range("K132").Value = 4458.6975
$ws.Range("L132").Value = 7352.678400000001
$ws.Range("M132").Value = -1928.6975
$ws.Range("N132").Value = -12412.6784

$ws = $wb.Worksheets("BSM")
$ws.Range("H56").Value = 36000
$ws.Range("I56").Value = 36000
$ws.Range("K56").Value = 36000
$ws.Range("M56").Value = -35261
$ws.Range("H86").Value = 134980.53
$ws.Range("I86").Value = 1517.3334
$ws.Range("J86").Value = 668833.3
$ws.Range("K86").Value = 1517.3334
$ws.Range("L86").Value = 668833.3
$ws.Range("M86").Value = -394.3334
$ws.Range("N86").Value = -671079.3
$ws.Range("H89").Value = 134980.53
$ws.Range("I89").Value = 1517.3334
$ws.Range("J89").Value = 668833.3
$ws.Range("K89").Value = 7586.666999999999
$ws.Range("L89").Value = 3344166.5
$ws.Range("M89").Value = -1970.666999999999
$ws.Range("N89").Value = -3355398.5
$ws.Range("H99").Value = 949
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 2795.5
$ws.Range("I105").Value = 2511.4167
$ws.Range("K105").Value = 2511.4167
$ws.Range("M105").Value = -764.4167000000002
$ws.Range("H107").Value = 3437.3635
$ws.Range("I107").Value = 3631.2
$ws.Range("J107").Value = 1499
$ws.Range("K107").Value = 3631.2
$ws.Range("L107").Value = 1499
$ws.Range("M107").Value = -1711.2
$ws.Range("N107").Value = -5339

$ws = $wb.Worksheets("CRP")
$ws.Range("H16").Value = 967.17645
$ws.Range("I16").Value = 896
$ws.Range("K16").Value = 896
$ws.Range("M16").Value = -609
$ws.Range("H20").Value = 49998.332
$ws.Range("J20").Value = 49998.332
$ws.Range("L20").Value = 49998.332
$ws.Range("N20").Value = -50470.332
$ws.Range("H30").Value = 49998.332
$ws.Range("J30").Value = 49998.332
$ws.Range("L30").Value = 49998.332
$ws.Range("N30").Value = -50180.332
$ws.Range("H31").Value = 2929.8096
$ws.Range("I31").Value = 2136.2727
$ws.Range("K31").Value = 2136.2727
$ws.Range("M31").Value = -1841.2727
$ws.Range("H34").Value = 2929.8096
$ws.Range("I34").Value = 2136.2727
$ws.Range("K34").Value = 2136.2727
$ws.Range("M34").Value = -1934.2727
$ws.Range("H99").Value = 3000
$ws.Range("I99").Value = 3500
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 3500
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -2002
$ws.Range("N99").Value = -4996
$ws.Range("H113").Value = 967.17645
$ws.Range("I113").Value = 896
$ws.Range("K113").Value = 896
$ws.Range("M113").Value = 1274
$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3500
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 10500
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -8030
$ws.Range("N126").Value = -10940
$ws.Range("H128").Value = 49998.332
$ws.Range("J128").Value = 49998.332
$ws.Range("L128").Value = 49998.332
$ws.Range("N128").Value = -59958.332
$ws.Range("H138").Value = 93795
$ws.Range("J138").Value = 93795
$ws.Range("L138").Value = 93795
$ws.Range("N138").Value = -104075

$ws = $wb.Worksheets("CUL")
$ws.Range("H68").Value = 973.125
$ws.Range("I68").Value = 797.5
$ws.Range("K68").Value = 2392.5
$ws.Range("M68").Value = -1581.5
$ws.Range("H71").Value = 973.125
$ws.Range("I71").Value = 797.5
$ws.Range("K71").Value = 7177.5
$ws.Range("M71").Value = -3121.5
$ws.Range("H122").Value = 1612.0625
$ws.Range("J122").Value = 1861
$ws.Range("L122").Value = 16749
$ws.Range("N122").Value = -21649
$ws.Range("H137").Value = 6216.1577
$ws.Range("I137").Value = 2943.8333
$ws.Range("J137").Value = 7726.4614
$ws.Range("K137").Value = 8831.499899999999
$ws.Range("L137").Value = 23179.3842
$ws.Range("M137").Value = -3731.499899999999
$ws.Range("N137").Value = -33379.3842

$ws = $wb.Worksheets("GSM")
$ws.Range("H113").Value = 1462.875
$ws.Range("I113").Value = 876
$ws.Range("K113").Value = 876
$ws.Range("M113").Value = 1294
$ws.Range("H122").Value = 1514.1904
$ws.Range("I122").Value = 1229.0834
$ws.Range("K122").Value = 3687.2502
$ws.Range("M122").Value = -1237.2502

$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 3189.6924
$ws.Range("I7").Value = 3267.3333
$ws.Range("K7").Value = 3267.3333
$ws.Range("M7").Value = -3155.3333
$ws.Range("H61").Value = 2141.45
$ws.Range("I61").Value = 1962.7222
$ws.Range("K61").Value = 1962.7222
$ws.Range("M61").Value = -1760.7222
$ws.Range("H68").Value = 2520.7646
$ws.Range("I68").Value = 1911.7693
$ws.Range("K68").Value = 1911.7693
$ws.Range("M68").Value = -1162.7693
$ws.Range("H71").Value = 2520.7646
$ws.Range("I71").Value = 1911.7693
$ws.Range("K71").Value = 9558.8465
$ws.Range("M71").Value = -5814.8465
$ws.Range("I93").Value = 1047
$ws.Range("K93").Value = 1047
$ws.Range("M93").Value = 201
$ws.Range("H113").Value = 2141.45
$ws.Range("I113").Value = 1962.7222
$ws.Range("K113").Value = 1962.7222
$ws.Range("M113").Value = 207.2778000000001
$ws.Range("H126").Value = 3189.6924
$ws.Range("I126").Value = 3267.3333
$ws.Range("K126").Value = 9801.999899999999
$ws.Range("M126").Value = -7331.999899999999
$ws.Range("H132").Value = 3363.2666
$ws.Range("I132").Value = 3007.739
$ws.Range("K132").Value = 9023.217000000001
$ws.Range("M132").Value = -6493.217000000001

$ws = $wb.Worksheets("WVR")
$ws.Range("H7").Value = 80005
$ws.Range("J7").Value = 80005
$ws.Range("L7").Value = 80005
$ws.Range("N7").Value = -80231
$ws.Range("H122").Value = 33453.152
$ws.Range("I122").Value = 41614.2
$ws.Range("K122").Value = 124842.6
$ws.Range("M122").Value = -122392.6
$ws.Range("H126").Value = 6964
$ws.Range("I126").Value = 6964
$ws.Range("K126").Value = 20892
$ws.Range("M126").Value = -18422
